$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.088.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.52%  "
$ws.Range("D3").Value = "'2.994.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'580.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "'162.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.64%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +3.31%  "
$ws.Range("D9").Value = "'2.991.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("D10").Value = "'6.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.85%  "
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("E12").Value = "  +5.30%  "
$ws.Range("D13").Value = "'0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.70%  "
$ws.Range("D14").Value = "'34.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.15%  "
$ws.Range("E15").Value = "  -0.72%  "
$ws.Range("D16").Value = "'66.123.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.58%  "
$ws.Range("D17").Value = "'3.490.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.16%  "
$ws.Range("E18").Value = "  +4.25%  "
$ws.Range("D19").Value = "'2.997.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.50%  "
$ws.Range("D20").Value = "'452.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.04%  "
$ws.Range("E21").Value = "  +5.74%  "
$ws.Range("E22").Value = "  +4.05%  "
$ws.Range("E23").Value = "  +6.64%  "
$ws.Range("D24").Value = "'82.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("D25").Value = "'2.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.62%  "
$ws.Range("D26").Value = "'12.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.86%  "
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'8.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +13.70%  "
$ws.Range("E30").Value = "  +18.58%  "
$ws.Range("E31").Value = "  +5.00%  "
$ws.Range("D32").Value = "'0.0000104"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.97%  "
$ws.Range("D33").Value = "'27.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.51%  "
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'0.990"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("D37").Value = "'5.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.76%  "
$ws.Range("D38").Value = "'2.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.74%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("E41").Value = "  +15.90%  "
$ws.Range("D42").Value = "'44.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.54%  "
$ws.Range("E43").Value = "  +6.69%  "
$ws.Range("D44").Value = "'8.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.51%  "
$ws.Range("D45").Value = "'397.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.70%  "
$ws.Range("E46").Value = "  +5.68%  "
$ws.Range("D47").Value = "'2.769.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").Value = "'133.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D50").Value = "'23.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.24%  "
$ws.Range("E51").Value = "  +3.72%  "
Write-Output "Applied crypto price/volume updates"
